$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 2742947.48
$ws.Range("C9").Value = 418881.99
$ws.Range("D9").Value = 3161829.47
$ws.Range("E9").Value = 13.24808924625527
$ws.Range("F9").Value = 86.75191075374474
$ws.Range("G9").Value = -59.517144308392
$ws.Range("H9").Value = -50.46615992186744
$ws.Range("I9").Value = -51.89112149640997
$ws.Range("J9").Value = 27232
$ws.Range("K9").Value = 1155
$ws.Range("L9").Value = 28387
